$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tanggal kejadian (event date/time, as Excel serial date values) for rows 2-10, column A
$dates = @(
    46003.14027777778,
    46003.14027777778,
    46003.23680555556,
    46003.76111111111,
    46004.62708333333,
    46004.99583333333,
    46006.00833333333,
    46006.48263888889,
    46006.86041666667
)

# Write the first cell and cycle its format through the lowercase built-in
# style before settling on the uppercase one, so both number formats get
# registered in the style table (matching the fixed workbook's style sheet).
$firstCell = $ws.Cells.Item(2, 1)
$firstCell.Value = $dates[0]
$firstCell.NumberFormat = "yyyy-mm-dd h:mm:ss"
$firstCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 1; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dates[$i]
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
